# Update BunkerPrices at 2025-03-19 14:59
#
# 1. Swap the "New York" / "Montevideo" columns (Y and Z) - header + data rows 2-4
# 2. Change AH4's number format from date-only (style 3) to date-time (style 2)
# 3. Append a new data row (row 5) for 2025-03-13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap columns Y (25) and Z (26) -------------------------------------

# Header row
$ws.Cells.Item(1, 25).Value = "Montevideo"
$ws.Cells.Item(1, 26).Value = "New York"

# Data rows 2-4: swap the Y/Z values
for ($r = 2; $r -le 4; $r++) {
    $yVal = $ws.Cells.Item($r, 25).Value()
    $zVal = $ws.Cells.Item($r, 26).Value()
    $ws.Cells.Item($r, 25).Value = $zVal
    $ws.Cells.Item($r, 26).Value = $yVal
}

# --- 2. AH4 (col 34) number format: YYYY-MM-DD -> YYYY-MM-DD HH:MM:SS ------

$ws.Cells.Item(4, 34).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 3. Append new row 5 ----------------------------------------------------

$row5 = @(577, 644, 502, 578, 620, 643, 509, 520, 570, 524, 585, 514, 531, 881, 578, 530, 509, 529, 608, 649, 582, 490, 555, 530, 551, 529, 509, 550, 578.5, 520, 515, 531, 490, 45729, 511, 557, 524, 765, 644, 613, 503, 640, 775, 518, 502, 565, 576, 635)

for ($i = 0; $i -lt $row5.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(5, $col).Value = $row5[$i]
}

# AH5 (col 34) keeps the plain-date format (style 3), matching AH2/AH3's previous sibling AH4
$ws.Cells.Item(5, 34).NumberFormat = "YYYY-MM-DD"
